$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Switch the field-name convention in both sheets' JETT template cells
#    from the mocked-up camelCase bean getters (employee.firstName, etc.)
#    to the HSQLDB-backed snake_case column names (employee.first_name, etc.)
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("Query")
$ws2 = $wb.Worksheets.Item("Prepared")

$ws2.Range("A3").Value = '<jt:forEach items="${jdbc.execQuery(''SELECT * FROM employee WHERE title = ?'', titleSearch)}" var="employee" >${employee.first_name}'
$ws1.Range("B2").Value = '${employee.last_name}'
$ws1.Range("A2").Value = '<jt:forEach items="${jdbc.execQuery(''SELECT * FROM employee'')}" var="employee" >${employee.first_name}'
$ws1.Range("F2").Value = '${employee.catch_phrase}'
$ws1.Range("G2").Value = '${employee.is_a_manager}</jt:forEach>'
$ws2.Range("G3").Value = '${employee.is_a_manager}</jt:forEach></jt:forEach>'
$ws2.Range("B3").Value = '${employee.last_name}'
$ws2.Range("F3").Value = '${employee.catch_phrase}'

# ---------------------------------------------------------------------------
# 2. Move the active tab from "Prepared" back to "Query".
# ---------------------------------------------------------------------------

$ws1.Activate()
